# The document has a single 5-column table. Each block of data occupies one
# table row followed by three blank spacer rows, so the five data rows sit
# at table-row indices 1, 5, 9, 13, 17 (1-based). We address each cell
# directly by (row, column) and overwrite its text, which keeps every
# other run/paragraph property (fonts, size, justification, etc.) intact.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Data row 1 (table row 1)
$t.Cell(1,1).Range.Text = "79÷3=26, 1"
$t.Cell(1,2).Range.Text = "90÷5=18, 0"
$t.Cell(1,3).Range.Text = "58÷7=8, 2"
$t.Cell(1,4).Range.Text = "23÷8=2, 7"
$t.Cell(1,5).Range.Text = "33÷3=11, 0"

# Data row 2 (table row 5)
$t.Cell(5,1).Range.Text = "72÷9=8, 0"
$t.Cell(5,2).Range.Text = "16÷8=2, 0"
$t.Cell(5,3).Range.Text = "99÷8=12, 3"
$t.Cell(5,4).Range.Text = "61÷7=8, 5"
$t.Cell(5,5).Range.Text = "71÷6=11, 5"

# Data row 3 (table row 9)
$t.Cell(9,1).Range.Text = "96÷7=13, 5"
$t.Cell(9,2).Range.Text = "67÷4=16, 3"
$t.Cell(9,3).Range.Text = "49÷5=9, 4"
$t.Cell(9,4).Range.Text = "11÷3=3, 2"
$t.Cell(9,5).Range.Text = "77÷7=11, 0"

# Data row 4 (table row 13)
$t.Cell(13,1).Range.Text = "19÷3=6, 1"
$t.Cell(13,2).Range.Text = "35÷8=4, 3"
$t.Cell(13,3).Range.Text = "87÷7=12, 3"
$t.Cell(13,4).Range.Text = "19÷7=2, 5"
$t.Cell(13,5).Range.Text = "61÷2=30, 1"

# Data row 5 (table row 17)
$t.Cell(17,1).Range.Text = "83÷3=27, 2"
$t.Cell(17,2).Range.Text = "54÷2=27, 0"
$t.Cell(17,3).Range.Text = "84÷6=14, 0"
$t.Cell(17,4).Range.Text = "19÷5=3, 4"
$t.Cell(17,5).Range.Text = "84÷4=21, 0"
